$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fromCSV")

# Find the last used row/column on the sheet (data currently spans A1:V130)
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# 1) Update the shared "short-url" value (column B) across all data rows:
#    "8DWxu2" -> "3Nljol"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Text -eq "8DWxu2") {
        $cell.Value = "3Nljol"
    }
}

# 2) Update the refugees figure for Afghanistan / Turkmenistan / 2024
#    (row 128, column N "refugees"): 3509 -> 3409
$ws.Cells.Item(128, 14).Value = 3409

# 3) Remove the last two rows (129 and 130), which are no longer part of the
#    export. Delete from the bottom up so row numbers above aren't disturbed.
if ($lastRow -ge 130) {
    $ws.Rows.Item(130).Delete()
}
if ($lastRow -ge 129) {
    $ws.Rows.Item(129).Delete()
}
